$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 7 - "Lab #8: Register write port"
# Shape 3 = TextBox 37 (id=38), the "Reminders:" bullet list.
# ---------------------------------------------------------------------------
$s7   = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(3)
$tr7  = $shp7.TextFrame.TextRange

# Apply edits from the end of the text backwards so earlier offsets stay valid.

# "The input to the selection pin (rt) will determine..." -> "(rd)"
$tr7.Characters(203, 2).Text = "rd"

# "WE = write enable; set to 1 when  you want to write"
#   -> "WE = write enable; set to 1 when you want to write"   (collapse double space)
$tr7.Characters(63, 51).Text = "WE = write enable; set to 1 when you want to write"

# ---------------------------------------------------------------------------
# Slide 8 - "Lab #8: Making adjustments to register 7 (r7)" (incr7 slide)
# Shape 2 = TextBox 36 (id=37).
# ---------------------------------------------------------------------------
$s8   = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(2)
$tr8  = $shp8.TextFrame.TextRange
$origHeight8 = $shp8.Height

# "the data input should be ignored" -> "the data input to r7 should be ignored"
$tr8.Characters(354, 32).Text = "the data input to r7 should be ignored"

# The autofit textbox grows when the text gets longer; restore the original
# height so the shape geometry matches the source edit exactly.
$shp8.Height = $origHeight8

# ---------------------------------------------------------------------------
# Slide 9 - "Lab #8: Making adjustments to register 7 (r7)" (decr7 slide)
# Shape 2 = TextBox 36 (id=37).
# ---------------------------------------------------------------------------
$s9   = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(2)
$tr9  = $shp9.TextFrame.TextRange
$origHeight9 = $shp9.Height

# Apply edits from the end of the text backwards so earlier offsets stay valid.

# "...or the register's content minus 1 (r7-1)..." -> "...contents minus 1..."
$tr9.Characters(364, 195).Text = "For your circuit, this implies that you will need to add logic to determine whether the current register" + [char]8217 + "s contents (r7) or the register" + [char]8217 + "s contents minus 1 (r7-1) should be sent to the read ports."

# "the data input should be ignored" -> "the data input to r7 should be ignored"
$tr9.Characters(178, 32).Text = "the data input to r7 should be ignored"

# The autofit textbox grows when the text gets longer; restore the original
# height so the shape geometry matches the source edit exactly.
$shp9.Height = $origHeight9

# ---------------------------------------------------------------------------
# Slide 10 - "Lab #8: Making adjustments to register 7 (r7)" (testing slide)
# Shape 2 = TextBox 36 (id=37).
# ---------------------------------------------------------------------------
$s10   = $p.Slides.Item(10)
$shp10 = $s10.Shapes.Item(2)
$tr10  = $shp10.TextFrame.TextRange

# Apply edits from the end of the text backwards so earlier offsets stay valid.

# " r7 should auto-decrement, even with write enable off"
#   -> " r7 should auto-decrement, even when the write enable pin = 0 "
$tr10.Characters(180, 53).Text = " r7 should auto-decrement, even when the write enable pin = 0 "

# "The contents of r7 should increase even with write enable off."
#   -> "The contents of r7 should increase even when the write enable pin = 0"
$tr10.Characters(82, 62).Text = "The contents of r7 should increase even when the write enable pin = 0"
